$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) is formatted as Text so values like "1.830" or
# "29.482.03" are not auto-converted/truncated into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$rows = @(
    @{ Row=2; B='Bitcoin'; C='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D='29.482.03'; E='  +0.81%  ' },
    @{ Row=3; B='Ethereum'; C='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D='1.970.84'; E='  +3.65%  ' },
    @{ Row=4; B='TetherUSD'; C='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D='1.004'; E='  +0.16%  ' },
    @{ Row=5; B='BNB'; C='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D='326.76'; E='  +0.24%  ' },
    @{ Row=6; B='USDC'; C='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D='1.004'; E='  +0.21%  ' },
    @{ Row=7; B='XRP'; C='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D='0.4668'; E='  +0.64%  ' },
    @{ Row=8; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.3921'; E='  +0.11%  ' },
    @{ Row=9; B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='46.26'; E='  -0.86%  ' },
    @{ Row=10; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.07936'; E='  +0.59%  ' },
    @{ Row=11; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='0.9892'; E='  +0.06%  ' },
    @{ Row=12; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='22.78'; E='  +4.46%  ' },
    @{ Row=13; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='2.036.99'; E='  +7.52%  ' },
    @{ Row=14; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='7.189'; E='  +1.63%  ' },
    @{ Row=15; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='5.806'; E='  +1.18%  ' },
    @{ Row=16; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.07123'; E='  +1.87%  ' },
    @{ Row=17; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='87.81'; E='  -0.51%  ' },
    @{ Row=18; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.005'; E='  +0.23%  ' },
    @{ Row=19; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.000009929'; E='  -0.46%  ' },
    @{ Row=20; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='17.26'; E='  +0.89%  ' },
    @{ Row=21; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.005'; E='  +0.38%  ' },
    @{ Row=22; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='29.518.83'; E='  +0.94%  ' },
    @{ Row=23; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='5.547'; E='  +4.74%  ' },
    @{ Row=24; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='11.15'; E='  +0.74%  ' },
    @{ Row=25; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='2.229.41'; E='  +4.85%  ' },
    @{ Row=26; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='2.108'; E='  +0.39%  ' },
    @{ Row=27; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='158.56'; E='  +1.63%  ' },
    @{ Row=28; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='19.54'; E='  +0.54%  ' },
    @{ Row=29; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='5.770'; E='  -3.65%  ' },
    @{ Row=30; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='119.86'; E='  +1.01%  ' },
    @{ Row=31; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='1.907'; E='  +1.25%  ' },
    @{ Row=32; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.09407'; E='  +0.64%  ' },
    @{ Row=33; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.8941'; E='  -0.74%  ' },
    @{ Row=34; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='5.240'; E='  -0.32%  ' },
    @{ Row=35; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.322'; E='  -0.09%  ' },
    @{ Row=36; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='3.172'; E='  -1.28%  ' },
    @{ Row=37; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.05824'; E='  +0.79%  ' },
    @{ Row=38; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='1.169'; E='  -1.56%  ' },
    @{ Row=39; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.02106'; E='  +0.91%  ' },
    @{ Row=40; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='7.745'; E='  +0.54%  ' },
    @{ Row=41; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.5711'; E='  +0.07%  ' },
    @{ Row=42; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.1795'; E='  +0.55%  ' },
    @{ Row=43; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='9.660'; E='  -0.47%  ' },
    @{ Row=44; B='PEPE'; C='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; D='0.000002917'; E='  +40.69%  ' },
    @{ Row=45; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.760'; E='  +7.50%  ' },
    @{ Row=46; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='11.75'; E='  -1.24%  ' },
    @{ Row=47; B='Decentraland'; C='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D='0.5338'; E='  -0.31%  ' },
    @{ Row=48; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='2.160'; E='  -0.60%  ' },
    @{ Row=49; B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.06923'; E='  -1.46%  ' },
    @{ Row=50; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='1.830'; E='  -1.09%  ' },
    @{ Row=51; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='113.70'; E='  +0.50%  ' }
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
}
